# Timetable Wordle.xlsx - update log with the "14.01.26" follow-up entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 6 ("12.01.26"): time invested corrected from 7h to 9h.
$ws.Range("C6").Value = 9

# Row 9 ("14.01.26"): fill in the content + hours that were left blank.
$ws.Range("B9").Value = "Worked on fixing a logic error while solving wordle where guess with equal letters would all be colored even if mystery word didn" + [char]0x2019 + "t have that many letters. Also did research on GUI (specifically JavaFX)"
$ws.Range("C9").Value = 8

# Update the on-screen selection to match the saved view (cell F6 selected).
$ws.Activate()
$ws.Range("F6").Select()
